$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = "[Luciano%Gattinoni%gattinoniluciano@gmail.com%0, Davide%Chiumello%NULL%3, Sandra%Rossi%NULL%0]"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "BioMed Central"

$ws.Range("E3").Value = "[Luciano%Gattinoni%NULL%0, Silvia%Coppola%NULL%2, Silvia%Coppola%NULL%0, Massimo%Cressoni%NULL%1, Mattia%Busana%NULL%2, Mattia%Busana%NULL%0, Sandra%Rossi%NULL%0, Sandra%Rossi%NULL%0, Davide%Chiumello%NULL%0, Davide%Chiumello%NULL%0]"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = "American Thoracic Society"

$abstract = "Aerosol generating procedures (AGPs) may expose health care workers (HCWs) to pathogens causing acute respiratory infections (ARIs), but the risk of transmission of ARIs from AGPs is not fully known.`n We sought to determine the clinical evidence for the risk of transmission of ARIs to HCWs caring for patients undergoing AGPs compared with the risk of transmission to HCWs caring for patients not undergoing AGPs.`n We searched PubMed, EMBASE, MEDLINE, CINAHL, the Cochrane Library, University of York CRD databases, EuroScan, LILACS, Indian Medlars, Index Medicus for SE Asia, international health technology agencies and the Internet in all languages for articles from 01/01/1990 to 22/10/2010. Independent reviewers screened abstracts using pre-defined criteria, obtained full-text articles, selected relevant studies, and abstracted data.`n Disagreements were resolved by consensus.`n The outcome of interest was risk of ARI transmission.`n The quality of evidence was rated using the GRADE system.`n We identified 5 case-control and 5 retrospective cohort studies which evaluated transmission of SARS to HCWs.`n Procedures reported to present an increased risk of transmission included [n; pooled OR(95%CI)] tracheal intubation [n = 4 cohort; 6.6 (2.3, 18.9), and n = 4 case-control; 6.6 (4.1, 10.6)], non-invasive ventilation [n = 2 cohort; OR 3.1(1.4, 6.8)], tracheotomy [n = 1 case-control; 4.2 (1.5, 11.5)] and manual ventilation before intubation [n = 1 cohort; OR 2.8 (1.3, 6.4)].`n Other intubation associated procedures, endotracheal aspiration, suction of body fluids, bronchoscopy, nebulizer treatment, administration of O2, high flow O2, manipulation of O2 mask or BiPAP mask, defibrillation, chest compressions, insertion of nasogastric tube, and collection of sputum were not significant.`n Our findings suggest that some procedures potentially capable of generating aerosols have been associated with increased risk of SARS transmission to HCWs or were a risk factor for transmission, with the most consistent association across multiple studies identified with tracheal intubation.`n`n"

$ws.Range("D4").Value = $abstract
$ws.Range("E4").Value = "[Khai%Tran%NULL%0, Karen%Cimon%NULL%1, Melissa%Severn%NULL%1, Carmem L.%Pessoa-Silva%NULL%1, John%Conly%NULL%1, Malcolm Gracie%Semple%NULL%2, Malcolm Gracie%Semple%NULL%0]"
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = "Public Library of Science"
